$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 3553
$ws.Range("J97").Value = 3553
$ws.Range("L97").Value = 10659
$ws.Range("N97").Value = -11651

$ws.Range("H127").Value = 1880.7368
$ws.Range("I127").Value = 912.1429
$ws.Range("J127").Value = 2445.75
$ws.Range("K127").Value = 2736.4287
$ws.Range("L127").Value = 7337.25
$ws.Range("M127").Value = 2223.5713
$ws.Range("N127").Value = -17257.25

$ws.Range("H129").Value = 928
$ws.Range("J129").Value = 999.1667
$ws.Range("L129").Value = 2997.5001
$ws.Range("N129").Value = -12997.5001

$ws.Range("H138").Value = 3651.3289
$ws.Range("I138").Value = 2384.0908
$ws.Range("J138").Value = 3865.7847
$ws.Range("K138").Value = 7152.2724
$ws.Range("L138").Value = 11597.3541
$ws.Range("M138").Value = -2012.2724
$ws.Range("N138").Value = -21877.3541

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 32259430
$ws.Range("I61").Value = 35715440
$ws.Range("J61").Value = 3342.6667
$ws.Range("K61").Value = 35715440
$ws.Range("L61").Value = 3342.6667
$ws.Range("M61").Value = -35715228
$ws.Range("N61").Value = -3766.6667

$ws.Range("H97").Value = 9755.182
$ws.Range("I97").Value = 756.3333
$ws.Range("J97").Value = 50250
$ws.Range("K97").Value = 756.3333
$ws.Range("L97").Value = 50250
$ws.Range("M97").Value = -260.3333
$ws.Range("N97").Value = -51242

$ws.Range("H136").Value = 32259430
$ws.Range("I136").Value = 35715440
$ws.Range("J136").Value = 3342.6667
$ws.Range("K136").Value = 107146320
$ws.Range("L136").Value = 10028.0001
$ws.Range("M136").Value = -107143770
$ws.Range("N136").Value = -15128.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 19231370
$ws.Range("I94").Value = 19231370
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 19231370
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -19230919
$ws.Range("N94").ClearContents()  # was -1161.33334

$ws.Range("H134").Value = 3069.2156
$ws.Range("I134").Value = 839.5625
$ws.Range("J134").Value = 6824.421
$ws.Range("K134").Value = 2518.6875
$ws.Range("L134").Value = 20473.263
$ws.Range("M134").Value = 16.3125
$ws.Range("N134").Value = -25543.263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()  # was -1176

$ws.Range("H7").Value = 287.7143
$ws.Range("I7").Value = 114.75
$ws.Range("J7").Value = 518.3333
$ws.Range("K7").Value = 114.75
$ws.Range("L7").Value = 518.3333
$ws.Range("M7").Value = -1.75
$ws.Range("N7").Value = -744.3333

$ws.Range("H20").Value = 48671.6
$ws.Range("J20").Value = 48671.6
$ws.Range("L20").Value = 48671.6
$ws.Range("N20").Value = -49143.6

$ws.Range("H30").Value = 48671.6
$ws.Range("J30").Value = 48671.6
$ws.Range("L30").Value = 48671.6
$ws.Range("N30").Value = -48853.6

$ws.Range("H31").Value = 1433.7297
$ws.Range("I31").Value = 1433.7297
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1433.7297
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1138.7297
$ws.Range("N31").ClearContents()  # was -4338.3333

$ws.Range("H34").Value = 1433.7297
$ws.Range("I34").Value = 1433.7297
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1433.7297
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1231.7297
$ws.Range("N34").ClearContents()  # was -4152.3333

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()  # was -26250

$ws.Range("H59").Value = 29760
$ws.Range("J59").Value = 29760
$ws.Range("L59").Value = 29760
$ws.Range("N59").Value = -32050

$ws.Range("H60").Value = 25000
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022

$ws.Range("H68").Value = 22750
$ws.Range("J68").Value = 22750
$ws.Range("L68").Value = 22750
$ws.Range("N68").Value = -24248

$ws.Range("H71").Value = 22750
$ws.Range("J71").Value = 22750
$ws.Range("L71").Value = 68250
$ws.Range("N71").Value = -75738

$ws.Range("H128").Value = 48671.6
$ws.Range("J128").Value = 48671.6
$ws.Range("L128").Value = 48671.6
$ws.Range("N128").Value = -58631.6

$ws.Range("H141").Value = 533036.3
$ws.Range("J141").Value = 533036.3
$ws.Range("L141").Value = 533036.3
$ws.Range("N141").Value = -543396.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 440.26315
$ws.Range("I7").Value = 471.4
$ws.Range("J7").Value = 323.5
$ws.Range("K7").Value = 1414.2
$ws.Range("L7").Value = 970.5
$ws.Range("M7").Value = -1302.2
$ws.Range("N7").Value = -1194.5

$ws.Range("H34").Value = 1977.7273
$ws.Range("J34").Value = 2997.1428
$ws.Range("L34").Value = 8991.4284
$ws.Range("N34").Value = -9159.4284

$ws.Range("H113").Value = 711.26086
$ws.Range("I113").Value = 480
$ws.Range("J113").Value = 721.7727
$ws.Range("K113").Value = 1440
$ws.Range("L113").Value = 2165.3181
$ws.Range("M113").Value = 730
$ws.Range("N113").Value = -6505.3181

$ws.Range("H137").Value = 25869090
$ws.Range("I137").Value = 75002770
$ws.Range("J137").Value = 9261.368
$ws.Range("K137").Value = 225008310
$ws.Range("L137").Value = 27784.104
$ws.Range("M137").Value = -225003210
$ws.Range("N137").Value = -37984.104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5418.5884
$ws.Range("I132").Value = 6077.2
$ws.Range("J132").Value = 3589.111
$ws.Range("K132").Value = 18231.6
$ws.Range("L132").Value = 10767.333
$ws.Range("M132").Value = -15701.6
$ws.Range("N132").Value = -15827.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1232.5
$ws.Range("J46").Value = 1700
$ws.Range("L46").Value = 1700
$ws.Range("N46").Value = -2076

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()  # was -43767.5

$ws.Range("H140").Value = 58118.285
$ws.Range("J140").Value = 58118.285
$ws.Range("L140").Value = 58118.285
$ws.Range("N140").Value = -68478.285

$ws.Range("H141").Value = 47550.418
$ws.Range("J141").Value = 46418.637
$ws.Range("L141").Value = 46418.637
$ws.Range("N141").Value = -56778.637

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1842.9048
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 4000
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 1842.9048
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 20000
$ws.Range("N84").Value = -30608

$ws.Range("H136").Value = 1375.2858
$ws.Range("I136").Value = 679.1429
$ws.Range("K136").Value = 2037.4287
$ws.Range("M136").Value = 512.5712999999998
